$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Insert a new row at position 4; existing rows 4-7 shift down to 5-8
$ws.Rows.Item(4).Insert()

# Refresh dimension / cell values for rows 2-8 per latest scrape (2026-01-06 12:40:20 JST)
# Row 2
$ws.Range("A2").Value = "2026-01-06 12:40:20"
$ws.Range("B2").Value = "法人向け生成AIサービス(RAG・議事録機能)の設計・開発を支援エンジニア募集(AI/バックエンド)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5445159"
$ws.Range("G2").Value = 368
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# Row 3
$ws.Range("A3").Value = "2026-01-06 12:40:20"
$ws.Range("B3").Value = "B2B向け生成AIサービス(チャット・RAG)の新規開発プロジェクト推進を支援してくださるPM募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5445154"
$ws.Range("G3").Value = 368
$ws.Range("H3").Value = "🔥AI,Ai ◆開発"

# Row 4
$ws.Range("A4").Value = "2026-01-06 12:40:20"
$ws.Range("B4").Value = "​【1万〜3万円/BASE経験者】アリエク・ネッシー等のCSVをBASE用に変換・加工するツール作成"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5465992"
$ws.Range("G4").Value = 65
$ws.Range("H4").Value = "◆ツール"

# Row 5
$ws.Range("A5").Value = "2026-01-06 12:40:20"
$ws.Range("B5").Value = "初回 【急募】ECサイトの要件定義や基本設計ができる方を募集(1人月、フルリモート可、2025年12月〜)"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5425629"
$ws.Range("G5").Value = 45
$ws.Range("H5").Value = "◇サイト"

# Row 6
$ws.Range("A6").Value = "2026-01-06 12:40:20"
$ws.Range("B6").Value = "現行のシステムに追加要素"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5465878"
$ws.Range("G6").Value = 33

# Row 7
$ws.Range("A7").Value = "2026-01-06 12:40:20"
$ws.Range("B7").Value = "《長期レギュラー》公的機関Web運用の要となる、ディレクター募集"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5465685"
$ws.Range("G7").Value = 18

# Row 8
$ws.Range("A8").Value = "2026-01-06 12:40:20"
$ws.Range("B8").Value = "ActiveDirectoryの移行(フェーズ1)"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5465836"
$ws.Range("G8").Value = 13

# Clear H6/H7/H8 leftovers copied down by the row insert (source rows had no H value)
$ws.Range("H6").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("H8").ClearContents()

# Rebuild hyperlinks cleanly so relationship ids line up with the (now shifted) URL column
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5445159") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5445154") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5465992") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5425629") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5465878") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5465685") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5465836") | Out-Null
$ws.Range("F2:F8").Style = "Hyperlink"

Write-Output "done"
